# Insert a new data row above row 88 (Feria Lagunitas de Puerto Montt - Cilantro),
# shifting all subsequent rows down by one, and populate the new row with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 88; existing row 88..181 shift to 89..182
$ws.Rows.Item(88).Insert()

# Populate the newly inserted row 88 with the new record's data
$ws.Range("A88").Value2 = 4
$ws.Range("B88").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C88").Value2 = "Los Lagos"
$ws.Range("D88").Value2 = 44494
$ws.Range("E88").Value2 = 10
$ws.Range("F88").Value2 = 100112040
$ws.Range("G88").Value2 = "Cilantro"
$ws.Range("H88").Value2 = "Sin especificar"
$ws.Range("I88").Value2 = "Primera"
$ws.Range("J88").Value2 = 100
$ws.Range("K88").Value2 = 10000
$ws.Range("L88").Value2 = 10000
$ws.Range("M88").Value2 = 10000
$ws.Range("N88").Value2 = "$/caja 36 atados"
$ws.Range("O88").Value2 = "Región Metropolitana"
$ws.Range("P88").Value2 = 278
$ws.Range("Q88").Value2 = 36
$ws.Range("R88").Value2 = "Hortaliza"
